# The workbook originally contains two sheets:
#   "selectActionTest1" (first/inactive tab)
#   "selectActionTest"  (second tab, was the active one)
#
# This change removes the "selectActionTest1" sheet entirely, leaving
# "selectActionTest" as the sole (and therefore active/selected) sheet.
# It also clears the stray "Inpatient Ward" reference that was left in
# cell E6 of "selectActionTest" (that value is no longer used once the
# first sheet - the only other place that referenced the same shared
# string - is gone), and updates the remembered selection on that sheet
# from D5 to D6.

$wb = $excel.ActiveWorkbook

# Remove the obsolete "selectActionTest1" worksheet.
$ws1 = $wb.Worksheets.Item("selectActionTest1")
$ws1.Delete()

# The remaining worksheet becomes the only / active sheet.
$ws = $wb.Worksheets.Item("selectActionTest")
$ws.Activate()

# Clear the leftover "Inpatient Ward" value in E6 - no longer referenced.
$ws.Range("E6").ClearContents()

# Update the remembered selection to D6.
$ws.Range("D6").Select()
